# Update specific imputed values in the KNN result data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value  = 7.593000000000001
$ws.Range("B7").Value  = 5.024
$ws.Range("C7").Value  = -13.208
$ws.Range("C15").Value = -13.491
$ws.Range("B16").Value = 5.235
$ws.Range("C21").Value = -12.518
$ws.Range("C22").Value = -12.961
$ws.Range("C23").Value = -12.223
$ws.Range("B28").Value = 5.709000000000001
$ws.Range("B29").Value = 5.065
$ws.Range("B32").Value = 6.387999999999999
$ws.Range("C34").Value = -11.997
$ws.Range("B40").Value = 9.044999999999998
$ws.Range("C43").Value = -13.516
$ws.Range("C45").Value = -13.03
$ws.Range("C50").Value = -13.812
$ws.Range("C51").Value = -11.276
$ws.Range("B52").Value = 5.546
$ws.Range("B57").Value = 5.093000000000001
$ws.Range("B66").Value = 5.054
$ws.Range("C66").Value = -10.883
$ws.Range("C67").Value = -11.395
$ws.Range("C79").Value = -12.117
$ws.Range("C84").Value = -13.72
$ws.Range("C92").Value = -11.139
$ws.Range("C97").Value = -12.349
$ws.Range("B100").Value = 5.587999999999999
